$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Insert three new columns before column E. This shifts the
#    existing PET/PEESE/Robustness data (old E..R) right to H..U,
#    carrying values, shared-strings, styles and merged cells along.
# ---------------------------------------------------------------
$ws.Range("E1:G1").EntireColumn.Insert()

# ---------------------------------------------------------------
# 2. Populate the new "Craig" / "Reproduce" columns (E,F,G)
#    (values are entered in the order the original author typed
#    them so new shared-string indices line up: Craig, r+, Reproduce)
# ---------------------------------------------------------------

# Header row 1: merged E1:F1 = "Craig"
$ws.Range("E1:F1").Merge()
$ws.Range("E1").Value = "Craig"
$ws.Range("E1:F1").HorizontalAlignment = -4108   # xlCenter

# Header row 2: k, r+, r+
$ws.Range("E2").Value = "k"
$ws.Range("F2").Value = "r+"
$ws.Range("G2").Value = "r+"

# Header row 1: G1 = "Reproduce"
$ws.Range("G1").Value = "Reproduce"
$ws.Range("G1").HorizontalAlignment = -4108      # xlCenter

# Data rows 3-14 (row, k, r+, r+)
$data = @(
    , @(3,  21, 0.29399999999999998, 0.31)
    , @(4,  37, 0.18099999999999999, 0.19)
    , @(5,  11, 0.10100000000000001, 0.11)
    , @(6,  19, 0.11,                0.15)
    , @(7,  27, 0.21,                0.22)
    , @(8,  45, 0.18,                0.19)
    , @(9,  40, 0.26200000000000001, 0.26)
    , @(10, 81, 0.189,               0.18)
    , @(11, 24, 0.217,               0.22)
    , @(12, 48, 0.20699999999999999, 0.21)
    , @(13, 27, 0.183,               0.16)
    , @(14, 38, 0.16400000000000001, 0.15)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("E$r").Value = $row[1]
    $ws.Range("F$r").Value = $row[2]
    $ws.Range("G$r").Value = $row[3]
}

# ---------------------------------------------------------------
# 3. Append the two new reviewer-comment rows at the bottom
# ---------------------------------------------------------------
$ws.Range("A21").Value = "Getting different results in reproduction for AggCog cross-sectional b/c of exclusion of Matsuzaki."
$ws.Range("A22").Value = "Must be some experimental effect sizes still showing up on several rows for AggBeh"

# ---------------------------------------------------------------
# 4. Update the view: selection moved to A23
# ---------------------------------------------------------------
$ws.Range("A23").Select()
